$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure new date cells (D228, D229) use the same custom date format as existing date column cells,
# so Excel reuses the existing style index instead of minting a new one.
$ws.Range("D228:D229").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 211
$ws.Range("A211").Value = 11
$ws.Range("B211").Value = 'Vega Monumental Concepción'
$ws.Range("C211").Value = 'Bíobío'
$ws.Range("D211").Value = '1/7/2022'
$ws.Range("E211").Value = 8
$ws.Range("F211").Value = 100112002
$ws.Range("G211").Value = 'Pimiento'
$ws.Range("H211").Value = 'Zafiro rojo'
$ws.Range("I211").Value = 'Primera'
$ws.Range("J211").Value = 100
$ws.Range("K211").Value = 21000
$ws.Range("L211").Value = 22000
$ws.Range("M211").Value = 21500
$ws.Range("N211").Value = '$/caja 15 kilos'
$ws.Range("O211").Value = 'Región de Arica y Parinacota'
$ws.Range("P211").Value = 1433
$ws.Range("Q211").Value = 15
$ws.Range("R211").Value = 'Hortaliza'

# Row 212
$ws.Range("A212").Value = 11
$ws.Range("B212").Value = 'Vega Monumental Concepción'
$ws.Range("C212").Value = 'Bíobío'
$ws.Range("D212").Value = '1/7/2022'
$ws.Range("E212").Value = 8
$ws.Range("F212").Value = 100112002
$ws.Range("G212").Value = 'Pimiento'
$ws.Range("H212").Value = 'Zafiro verde'
$ws.Range("I212").Value = 'Primera'
$ws.Range("J212").Value = 100
$ws.Range("K212").Value = 16000
$ws.Range("L212").Value = 17000
$ws.Range("M212").Value = 16500
$ws.Range("N212").Value = '$/caja 15 kilos'
$ws.Range("O212").Value = 'Región de Arica y Parinacota'
$ws.Range("P212").Value = 1100
$ws.Range("Q212").Value = 15
$ws.Range("R212").Value = 'Hortaliza'

# Row 213
$ws.Range("A213").Value = 11
$ws.Range("B213").Value = 'Vega Monumental Concepción'
$ws.Range("C213").Value = 'Bíobío'
$ws.Range("D213").Value = '9/24/2021'
$ws.Range("E213").Value = 8
$ws.Range("F213").Value = 100112002
$ws.Range("G213").Value = 'Pimiento'
$ws.Range("H213").Value = 'Morrón rojo'
$ws.Range("I213").Value = 'Primera'
$ws.Range("J213").Value = 100
$ws.Range("K213").Value = 57000
$ws.Range("L213").Value = 58000
$ws.Range("M213").Value = 57500
$ws.Range("N213").Value = '$/caja 18 kilos'
$ws.Range("O213").Value = 'Provincia de Limarí'
$ws.Range("P213").Value = 3194
$ws.Range("Q213").Value = 18
$ws.Range("R213").Value = 'Hortaliza'

# Row 214
$ws.Range("A214").Value = 11
$ws.Range("B214").Value = 'Vega Monumental Concepción'
$ws.Range("C214").Value = 'Bíobío'
$ws.Range("D214").Value = '6/10/2021'
$ws.Range("E214").Value = 8
$ws.Range("F214").Value = 100112002
$ws.Range("G214").Value = 'Pimiento'
$ws.Range("H214").Value = 'Zafiro rojo'
$ws.Range("I214").Value = 'Primera'
$ws.Range("J214").Value = 200
$ws.Range("K214").Value = 15000
$ws.Range("L214").Value = 16000
$ws.Range("M214").Value = 15500
$ws.Range("N214").Value = '$/caja 15 kilos'
$ws.Range("O214").Value = 'Región de Arica y Parinacota'
$ws.Range("P214").Value = 1033
$ws.Range("Q214").Value = 15
$ws.Range("R214").Value = 'Hortaliza'

# Row 215
$ws.Range("A215").Value = 11
$ws.Range("B215").Value = 'Vega Monumental Concepción'
$ws.Range("C215").Value = 'Bíobío'
$ws.Range("D215").Value = '6/10/2021'
$ws.Range("E215").Value = 8
$ws.Range("F215").Value = 100112002
$ws.Range("G215").Value = 'Pimiento'
$ws.Range("H215").Value = 'Zafiro verde'
$ws.Range("I215").Value = 'Primera'
$ws.Range("J215").Value = 200
$ws.Range("K215").Value = 11000
$ws.Range("L215").Value = 12000
$ws.Range("M215").Value = 11500
$ws.Range("N215").Value = '$/caja 15 kilos'
$ws.Range("O215").Value = 'Región de Arica y Parinacota'
$ws.Range("P215").Value = 767
$ws.Range("Q215").Value = 15
$ws.Range("R215").Value = 'Hortaliza'

# Row 216
$ws.Range("A216").Value = 11
$ws.Range("B216").Value = 'Vega Monumental Concepción'
$ws.Range("C216").Value = 'Bíobío'
$ws.Range("D216").Value = '6/24/2021'
$ws.Range("E216").Value = 8
$ws.Range("F216").Value = 100112002
$ws.Range("G216").Value = 'Pimiento'
$ws.Range("H216").Value = 'Cuatro cascos rojo'
$ws.Range("I216").Value = 'Primera'
$ws.Range("J216").Value = 100
$ws.Range("K216").Value = 18000
$ws.Range("L216").Value = 20000
$ws.Range("M216").Value = 19000
$ws.Range("N216").Value = '$/caja 18 kilos'
$ws.Range("O216").Value = 'Provincia de Limarí'
$ws.Range("P216").Value = 1056
$ws.Range("Q216").Value = 18
$ws.Range("R216").Value = 'Hortaliza'

# Row 217
$ws.Range("A217").Value = 11
$ws.Range("B217").Value = 'Vega Monumental Concepción'
$ws.Range("C217").Value = 'Bíobío'
$ws.Range("D217").Value = '6/24/2021'
$ws.Range("E217").Value = 8
$ws.Range("F217").Value = 100112002
$ws.Range("G217").Value = 'Pimiento'
$ws.Range("H217").Value = 'Cuatro cascos verde'
$ws.Range("I217").Value = 'Primera'
$ws.Range("J217").Value = 100
$ws.Range("K217").Value = 15000
$ws.Range("L217").Value = 16000
$ws.Range("M217").Value = 15500
$ws.Range("N217").Value = '$/caja 18 kilos'
$ws.Range("O217").Value = 'Provincia de Limarí'
$ws.Range("P217").Value = 861
$ws.Range("Q217").Value = 18
$ws.Range("R217").Value = 'Hortaliza'

# Row 218
$ws.Range("A218").Value = 11
$ws.Range("B218").Value = 'Vega Monumental Concepción'
$ws.Range("C218").Value = 'Bíobío'
$ws.Range("D218").Value = '6/18/2021'
$ws.Range("E218").Value = 8
$ws.Range("F218").Value = 100112002
$ws.Range("G218").Value = 'Pimiento'
$ws.Range("H218").Value = 'Zafiro amarillo'
$ws.Range("I218").Value = 'Primera'
$ws.Range("J218").Value = 100
$ws.Range("K218").Value = 33000
$ws.Range("L218").Value = 35000
$ws.Range("M218").Value = 34000
$ws.Range("N218").Value = '$/caja 15 kilos'
$ws.Range("O218").Value = 'Región de Arica y Parinacota'
$ws.Range("P218").Value = 2267
$ws.Range("Q218").Value = 15
$ws.Range("R218").Value = 'Hortaliza'

# Row 219
$ws.Range("A219").Value = 11
$ws.Range("B219").Value = 'Vega Monumental Concepción'
$ws.Range("C219").Value = 'Bíobío'
$ws.Range("D219").Value = '6/18/2021'
$ws.Range("E219").Value = 8
$ws.Range("F219").Value = 100112002
$ws.Range("G219").Value = 'Pimiento'
$ws.Range("H219").Value = 'Zafiro verde'
$ws.Range("I219").Value = 'Primera'
$ws.Range("J219").Value = 100
$ws.Range("K219").Value = 14000
$ws.Range("L219").Value = 15000
$ws.Range("M219").Value = 14500
$ws.Range("N219").Value = '$/caja 15 kilos'
$ws.Range("O219").Value = 'Región de Arica y Parinacota'
$ws.Range("P219").Value = 967
$ws.Range("Q219").Value = 15
$ws.Range("R219").Value = 'Hortaliza'

# Row 220
$ws.Range("A220").Value = 11
$ws.Range("B220").Value = 'Vega Monumental Concepción'
$ws.Range("C220").Value = 'Bíobío'
$ws.Range("D220").Value = '12/29/2020'
$ws.Range("E220").Value = 8
$ws.Range("F220").Value = 100112002
$ws.Range("G220").Value = 'Pimiento'
$ws.Range("H220").Value = 'Zafiro rojo'
$ws.Range("I220").Value = 'Primera'
$ws.Range("J220").Value = 100
$ws.Range("K220").Value = 30000
$ws.Range("L220").Value = 32000
$ws.Range("M220").Value = 31000
$ws.Range("N220").Value = '$/caja 15 kilos'
$ws.Range("O220").Value = 'Región de Arica y Parinacota'
$ws.Range("P220").Value = 2067
$ws.Range("Q220").Value = 15
$ws.Range("R220").Value = 'Hortaliza'

# Row 221
$ws.Range("A221").Value = 11
$ws.Range("B221").Value = 'Vega Monumental Concepción'
$ws.Range("C221").Value = 'Bíobío'
$ws.Range("D221").Value = '12/29/2020'
$ws.Range("E221").Value = 8
$ws.Range("F221").Value = 100112002
$ws.Range("G221").Value = 'Pimiento'
$ws.Range("H221").Value = 'Zafiro verde'
$ws.Range("I221").Value = 'Primera'
$ws.Range("J221").Value = 100
$ws.Range("K221").Value = 23000
$ws.Range("L221").Value = 24000
$ws.Range("M221").Value = 23500
$ws.Range("N221").Value = '$/caja 15 kilos'
$ws.Range("O221").Value = 'Región de Arica y Parinacota'
$ws.Range("P221").Value = 1567
$ws.Range("Q221").Value = 15
$ws.Range("R221").Value = 'Hortaliza'

# Row 222
$ws.Range("A222").Value = 11
$ws.Range("B222").Value = 'Vega Monumental Concepción'
$ws.Range("C222").Value = 'Bíobío'
$ws.Range("D222").Value = '4/27/2021'
$ws.Range("E222").Value = 8
$ws.Range("F222").Value = 100112002
$ws.Range("G222").Value = 'Pimiento'
$ws.Range("H222").Value = 'Morrón rojo'
$ws.Range("I222").Value = 'Primera'
$ws.Range("J222").Value = 100
$ws.Range("K222").Value = 8000
$ws.Range("L222").Value = 9000
$ws.Range("M222").Value = 8500
$ws.Range("N222").Value = '$/caja 18 kilos'
$ws.Range("O222").Value = 'Provincia de Limarí'
$ws.Range("P222").Value = 472
$ws.Range("Q222").Value = 18
$ws.Range("R222").Value = 'Hortaliza'

# Row 223
$ws.Range("A223").Value = 11
$ws.Range("B223").Value = 'Vega Monumental Concepción'
$ws.Range("C223").Value = 'Bíobío'
$ws.Range("D223").Value = '4/27/2021'
$ws.Range("E223").Value = 8
$ws.Range("F223").Value = 100112002
$ws.Range("G223").Value = 'Pimiento'
$ws.Range("H223").Value = 'Zafiro rojo'
$ws.Range("I223").Value = 'Primera'
$ws.Range("J223").Value = 100
$ws.Range("K223").Value = 14000
$ws.Range("L223").Value = 15000
$ws.Range("M223").Value = 14500
$ws.Range("N223").Value = '$/caja 15 kilos'
$ws.Range("O223").Value = 'Región de Arica y Parinacota'
$ws.Range("P223").Value = 967
$ws.Range("Q223").Value = 15
$ws.Range("R223").Value = 'Hortaliza'

# Row 224
$ws.Range("A224").Value = 11
$ws.Range("B224").Value = 'Vega Monumental Concepción'
$ws.Range("C224").Value = 'Bíobío'
$ws.Range("D224").Value = '4/27/2021'
$ws.Range("E224").Value = 8
$ws.Range("F224").Value = 100112002
$ws.Range("G224").Value = 'Pimiento'
$ws.Range("H224").Value = 'Zafiro verde'
$ws.Range("I224").Value = 'Primera'
$ws.Range("J224").Value = 100
$ws.Range("K224").Value = 9000
$ws.Range("L224").Value = 10000
$ws.Range("M224").Value = 9500
$ws.Range("N224").Value = '$/caja 15 kilos'
$ws.Range("O224").Value = 'Región de Arica y Parinacota'
$ws.Range("P224").Value = 633
$ws.Range("Q224").Value = 15
$ws.Range("R224").Value = 'Hortaliza'

# Row 225
$ws.Range("A225").Value = 11
$ws.Range("B225").Value = 'Vega Monumental Concepción'
$ws.Range("C225").Value = 'Bíobío'
$ws.Range("D225").Value = '12/10/2021'
$ws.Range("E225").Value = 8
$ws.Range("F225").Value = 100112002
$ws.Range("G225").Value = 'Pimiento'
$ws.Range("H225").Value = 'Cuatro cascos verde'
$ws.Range("I225").Value = 'Primera'
$ws.Range("J225").Value = 180
$ws.Range("K225").Value = 11000
$ws.Range("L225").Value = 12000
$ws.Range("M225").Value = 11556
$ws.Range("N225").Value = '$/caja 18 kilos'
$ws.Range("O225").Value = 'Provincia de Limarí'
$ws.Range("P225").Value = 642
$ws.Range("Q225").Value = 18
$ws.Range("R225").Value = 'Hortaliza'

# Row 226
$ws.Range("A226").Value = 11
$ws.Range("B226").Value = 'Vega Monumental Concepción'
$ws.Range("C226").Value = 'Bíobío'
$ws.Range("D226").Value = '12/10/2021'
$ws.Range("E226").Value = 8
$ws.Range("F226").Value = 100112002
$ws.Range("G226").Value = 'Pimiento'
$ws.Range("H226").Value = 'Zafiro rojo'
$ws.Range("I226").Value = 'Primera'
$ws.Range("J226").Value = 180
$ws.Range("K226").Value = 25000
$ws.Range("L226").Value = 26000
$ws.Range("M226").Value = 25444
$ws.Range("N226").Value = '$/caja 18 kilos'
$ws.Range("O226").Value = 'Limache'
$ws.Range("P226").Value = 1414
$ws.Range("Q226").Value = 18
$ws.Range("R226").Value = 'Hortaliza'

# Row 227
$ws.Range("A227").Value = 11
$ws.Range("B227").Value = 'Vega Monumental Concepción'
$ws.Range("C227").Value = 'Bíobío'
$ws.Range("D227").Value = '3/17/2021'
$ws.Range("E227").Value = 8
$ws.Range("F227").Value = 100112002
$ws.Range("G227").Value = 'Pimiento'
$ws.Range("H227").Value = 'Cuatro cascos verde'
$ws.Range("I227").Value = 'Primera'
$ws.Range("J227").Value = 100
$ws.Range("K227").Value = 8000
$ws.Range("L227").Value = 9000
$ws.Range("M227").Value = 8500
$ws.Range("N227").Value = '$/caja 18 kilos'
$ws.Range("O227").Value = 'Provincia de Limarí'
$ws.Range("P227").Value = 472
$ws.Range("Q227").Value = 18
$ws.Range("R227").Value = 'Hortaliza'

# Row 228
$ws.Range("A228").Value = 11
$ws.Range("B228").Value = 'Vega Monumental Concepción'
$ws.Range("C228").Value = 'Bíobío'
$ws.Range("D228").Value = '3/17/2021'
$ws.Range("E228").Value = 8
$ws.Range("F228").Value = 100112002
$ws.Range("G228").Value = 'Pimiento'
$ws.Range("H228").Value = 'Morrón rojo'
$ws.Range("I228").Value = 'Primera'
$ws.Range("J228").Value = 100
$ws.Range("K228").Value = 8000
$ws.Range("L228").Value = 9000
$ws.Range("M228").Value = 8500
$ws.Range("N228").Value = '$/caja 18 kilos'
$ws.Range("O228").Value = 'Provincia de Limarí'
$ws.Range("P228").Value = 472
$ws.Range("Q228").Value = 18
$ws.Range("R228").Value = 'Hortaliza'

# Row 229
$ws.Range("A229").Value = 11
$ws.Range("B229").Value = 'Vega Monumental Concepción'
$ws.Range("C229").Value = 'Bíobío'
$ws.Range("D229").Value = '3/17/2021'
$ws.Range("E229").Value = 8
$ws.Range("F229").Value = 100112002
$ws.Range("G229").Value = 'Pimiento'
$ws.Range("H229").Value = 'Zafiro rojo'
$ws.Range("I229").Value = 'Primera'
$ws.Range("J229").Value = 100
$ws.Range("K229").Value = 14000
$ws.Range("L229").Value = 15000
$ws.Range("M229").Value = 14500
$ws.Range("N229").Value = '$/caja 15 kilos'
$ws.Range("O229").Value = 'Región de Arica y Parinacota'
$ws.Range("P229").Value = 967
$ws.Range("Q229").Value = 15
$ws.Range("R229").Value = 'Hortaliza'
